$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.573.57'
$ws.Range('E2').Value = '  +5.39%  '
$ws.Range('D3').Value = '2.301.52'
$ws.Range('E3').Value = '  +3.32%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.78'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.63'
$ws.Range('E6').Value = '  +12.09%  '
$ws.Range('E7').Value = '  +1.79%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('E9').Value = '  +5.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.57'
$ws.Range('E10').Value = '  +9.88%  '
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.45'
$ws.Range('E12').Value = '  +6.32%  '
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '2.649.19'
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('D15').Value = '2.298.48'
$ws.Range('E15').Value = '  +3.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.83'
$ws.Range('E16').Value = '  +3.09%  '
$ws.Range('E17').Value = '  +4.64%  '
$ws.Range('D18').Value = '46.591.19'
$ws.Range('E18').Value = '  +5.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.04'
$ws.Range('E19').Value = '  +6.13%  '
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').Value = '  +4.02%  '
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.18'
$ws.Range('E22').Value = '  +2.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.30'
$ws.Range('E23').Value = '  +6.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.90'
$ws.Range('E24').Value = '  +3.17%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.94'
$ws.Range('E26').Value = '  +4.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '42.45'
$ws.Range('E27').Value = '  +9.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.91'
$ws.Range('E29').Value = '  +5.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.07'
$ws.Range('E30').Value = '  +4.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.85'
$ws.Range('E31').Value = '  +14.50%  '
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '147.51'
$ws.Range('E33').Value = '  -2.57%  '
$ws.Range('E34').Value = '  +4.08%  '
$ws.Range('E35').Value = '  +14.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.115'
$ws.Range('E36').Value = '  +10.41%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.19'
$ws.Range('E38').Value = '  +19.71%  '
$ws.Range('E39').Value = '  +5.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.01'
$ws.Range('E40').Value = '  +11.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.36'
$ws.Range('E41').Value = '  +5.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0303'
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.00'
$ws.Range('E44').Value = '  +9.61%  '
$ws.Range('D45').Value = '1.812.90'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.30'
$ws.Range('E46').Value = '  +20.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.197'
$ws.Range('E47').Value = '  +5.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '73.44'
$ws.Range('E48').Value = '  +8.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.90'
$ws.Range('E49').Value = '  +6.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '95.84'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('D51').Value = '2.527.60'
$ws.Range('E51').Value = '  +3.38%  '
